$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet restructuring:
#    "Sheet2" (empty) becomes "email_solution" (filled with data).
#    A brand-new sheet "loginToForest" is inserted right before it.
# ---------------------------------------------------------------------------
$sheet2 = $wb.Worksheets.Item("Sheet2")
$loginSheet = $wb.Worksheets.Add($sheet2)
$loginSheet.Name = "loginToForest"

$emailSheet = $wb.Worksheets.Item("Sheet2")
$emailSheet.Name = "email_solution"

# ---------------------------------------------------------------------------
# 2. Cell values -- order chosen to reproduce the shared-string table order
# ---------------------------------------------------------------------------
$emailSheet.Range("A1").Value = "shortlink"
$emailSheet.Range("B1").Value = "forestlink"

$loginSheet.Range("A1").Value = "forestlink"
$loginSheet.Range("B1").Value = "runmode"
$loginSheet.Range("A2").Value = "https://admin.upnorway.com/39039/data/1571794/index/record/1571792/561/summary"
$loginSheet.Range("B2").Value = "Y"

$emailSheet.Range("A2").Value = "https://qa.upnorway.net/partner-experience-confirmation/e5af5490-daa6-11e9-9c60-6faa94e82755/message"
$emailSheet.Range("B2").Value = "https://admin.upnorway.com/39039/data/1571794/index/record/1571792/561/summary"

$emailSheet.Range("D1").Value = "reply"

$emailSheet.Range("C2").Value = "auto first message"
$emailSheet.Range("D2").Value = "auto first reply"
$emailSheet.Range("E2").Value = "Y"

$emailSheet.Range("C4").Value = "auto third message"
$emailSheet.Range("D4").Value = "auto third reply"

$emailSheet.Range("B3").Value = "https://admin.upnorway.com/39039/data/1571794/index/record/1571792/537/summary"
$emailSheet.Range("A3").Value = "https://qa.upnorway.net/partner-experience-confirmation/21ae5b60-c7df-11e9-8a1b-5da771ba9daa/message"
$emailSheet.Range("C3").Value = "arshad first message"
$emailSheet.Range("D3").Value = "arshad second reply 888"

# remaining cells re-use already-registered shared strings
$emailSheet.Range("C1").Value = "message"
$emailSheet.Range("E1").Value = "runmode"
$emailSheet.Range("E3").Value = "Y"
$emailSheet.Range("A4").Value = "https://qa.upnorway.net/partner-experience-confirmation/e5af5490-daa6-11e9-9c60-6faa94e82755/message"
$emailSheet.Range("B4").Value = "https://admin.upnorway.com/39039/data/1571794/index/record/1571792/561/summary"
$emailSheet.Range("E4").Value = "N"

# ---------------------------------------------------------------------------
# 3. Formatting: wrap text, row heights, column widths
# ---------------------------------------------------------------------------
$loginSheet.Range("A1").WrapText = $true
$loginSheet.Range("A2").WrapText = $true
$loginSheet.Rows(2).RowHeight = 144

$emailSheet.Range("A1:E1").WrapText = $true
$emailSheet.Range("A2:E4").WrapText = $true
$emailSheet.Rows(2).RowHeight = 72
$emailSheet.Rows(3).RowHeight = 72
$emailSheet.Rows(4).RowHeight = 72

$emailSheet.Columns("A").ColumnWidth = 25.166666666666668
$emailSheet.Columns("B:C").ColumnWidth = 33.608072916666664

# ---------------------------------------------------------------------------
# 4. Hyperlinks (href text equals the displayed URL text)
# ---------------------------------------------------------------------------
$loginSheet.Hyperlinks.Add($loginSheet.Range("A2"), "https://admin.upnorway.com/39039/data/1571794/index/record/1571792/561/summary") | Out-Null

$emailSheet.Hyperlinks.Add($emailSheet.Range("B2"), "https://admin.upnorway.com/39039/data/1571794/index/record/1571792/561/summary") | Out-Null
$emailSheet.Hyperlinks.Add($emailSheet.Range("A2"), "https://qa.upnorway.net/partner-experience-confirmation/e5af5490-daa6-11e9-9c60-6faa94e82755/message") | Out-Null
$emailSheet.Hyperlinks.Add($emailSheet.Range("B3"), "https://admin.upnorway.com/39039/data/1571794/index/record/1571792/537/summary") | Out-Null
$emailSheet.Hyperlinks.Add($emailSheet.Range("B4"), "https://admin.upnorway.com/39039/data/1571794/index/record/1571792/561/summary") | Out-Null
$emailSheet.Hyperlinks.Add($emailSheet.Range("A3"), "https://qa.upnorway.net/partner-experience-confirmation/21ae5b60-c7df-11e9-8a1b-5da771ba9daa/message") | Out-Null
$emailSheet.Hyperlinks.Add($emailSheet.Range("A4"), "https://qa.upnorway.net/partner-experience-confirmation/e5af5490-daa6-11e9-9c60-6faa94e82755/message") | Out-Null

# ---------------------------------------------------------------------------
# 5. Selections / active sheet
# ---------------------------------------------------------------------------
$loginSheet.Range("B2").Select()
$emailSheet.Activate()
$emailSheet.Range("F9").Select()
